$d = $word.ActiveDocument

# --- Fix "bug with empty notes": footnote 23 only contains a stray "།"
# (an effectively empty/meaningless note) together with a couple of stray
# space-runs around it. Remove the footnote (which also removes its
# in-text reference run) and then tidy up the leftover double space so the
# page-folio marker "[༡༠བ]" is separated from the poem by a single space,
# merged back into one run.

# The footnotes collection is 1-based and only contains "real" footnotes
# (21, 22, 23) in document order, so the last one is always the stray note
# we need to get rid of.
$fn = $d.Footnotes.Item($d.Footnotes.Count)
$fn.Delete()

# After the footnote reference is gone, the paragraph ends in:
#   ... ཀླུ་སྒྲུབ་ཀྱི་གླུ་རྫོགས་སོ།། །།<space><space>[༡༠བ]
# Collapse the doubled space down to one and fold everything back into a
# single run, matching how the rest of the paragraph is stored.
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute(
    "ངོ་། །ཀླུ་སྒྲུབ་ཀྱི་གླུ་རྫོགས་སོ།། །།  [༡༠བ]",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ངོ་། །ཀླུ་སྒྲུབ་ཀྱི་གླུ་རྫོགས་སོ།། །། [༡༠བ]", 2)
